$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the student names in column A with the new set of names,
# keeping course/grade/date columns untouched.
$ws.Range("A2").Value = "Ahmed Abdullah"
$ws.Range("A4").Value = "Ibrahim Mohamed"
$ws.Range("A3").Value = "Ahmed Ashraf"
$ws.Range("A5").Value = "Mohamed Omara"
$ws.Range("A6").Value = "Ahmed Ayman"

$ws.Range("N9").Select()
